$d = $word.ActiveDocument

# Locate the paragraph that holds the error-message run so we can rebuild
# it precisely: split "e query aql" into four runs, add the spell-check
# proofErr markers, drop the w:b on the trailing run, and shorten its text
# to "ecore".
$marker = "Couldn't find the self variable"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($marker)) {
        $targetIndex = $i
    }
}

$targetRange = $d.Paragraphs.Item($targetIndex).Range

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>Template de test pour les balises d</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">e </w:t></w:r>' +
  '<w:r><w:t>query</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>aql</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> : </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>ecore</w:t></w:r>' +
  '</w:p>'

# InsertXML replaces the target range's contents, so calling it on the whole
# paragraph range swaps the old paragraph for the rebuilt one in place.
[void]$targetRange.InsertXML($newParagraphXml)
